$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)

$rng3 = $ftr.Range.Duplicate
$found3 = $rng3.Find.Execute("Last update: ")
$start = $rng3.End

$rng2 = $ftr.Range.Duplicate
$found2 = $rng2.Find.Execute("support@bitvis.no")
$end = $rng2.Start

Write-Output "start=$start end=$end"

$dateRange = $ftr.Range.Duplicate
$dateRange.Start = $start
$dateRange.End = $end
Write-Output "dateRange.Text=[$($dateRange.Text)]"
